$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text storage (avoid Excel auto-converting numeric/percent-looking
# strings into actual numbers) by forcing Text format before assignment.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "312.42"
$ws.Range("E2").Value = "5.25%"
$ws.Range("E3").Value = "7.28%"
$ws.Range("D4").Value = "5.117"
$ws.Range("E4").Value = "1.52%"
$ws.Range("D5").Value = "0.08000"
$ws.Range("E5").Value = "5.97%"
$ws.Range("D6").Value = "4.500"
$ws.Range("E6").Value = "2.62%"
$ws.Range("D7").Value = "1.657"
$ws.Range("E7").Value = "3.68%"
$ws.Range("E8").Value = "16.94%"
$ws.Range("D9").Value = "0.1296"
$ws.Range("E9").Value = "6.65%"
$ws.Range("E10").Value = "3.92%"
$ws.Range("D11").Value = "0.09388"
$ws.Range("E11").Value = "4.64%"
$ws.Range("D12").Value = "0.04223"
$ws.Range("E12").Value = "7.12%"
$ws.Range("E13").Value = "-1.11%"
$ws.Range("D14").Value = "0.001305"
$ws.Range("E14").Value = "1.70%"
$ws.Range("D15").Value = "0.005837"
$ws.Range("E15").Value = "-3.06%"
$ws.Range("E17").Value = "1.26%"
$ws.Range("D18").Value = "2.401"
$ws.Range("E18").Value = "-0.82%"
$ws.Range("E19").Value = "1.92%"
$ws.Range("D20").Value = "8.055"
$ws.Range("E20").Value = "1.61%"
$ws.Range("D21").Value = "0.1370"
$ws.Range("E21").Value = "-3.48%"
$ws.Range("D23").Value = "0.04199"
$ws.Range("E23").Value = "3.45%"
$ws.Range("D24").Value = "0.001274"
$ws.Range("E24").Value = "0.74%"
$ws.Range("D25").Value = "0.004598"
$ws.Range("E25").Value = "15.63%"
$ws.Range("D26").Value = "0.0001340"
$ws.Range("E26").Value = "8.99%"
$ws.Range("D38").Value = "0.02652"
$ws.Range("E38").Value = "10.28%"
$ws.Range("D39").Value = "0.05402"
$ws.Range("E39").Value = "3.70%"
$ws.Range("D40").Value = "0.005628"
$ws.Range("E40").Value = "-12.06%"
$ws.Range("D41").Value = "0.007806"
$ws.Range("E41").Value = "0.44%"
$ws.Range("D42").Value = "0.1412"
$ws.Range("E42").Value = "6.33%"
$ws.Range("D43").Value = "0.007331"
$ws.Range("E43").Value = "-3.03%"
$ws.Range("D44").Value = "0.007883"
$ws.Range("E44").Value = "0.58%"
$ws.Range("D45").Value = "0.3115"
$ws.Range("E45").Value = "-3.08%"
$ws.Range("D46").Value = "0.00006729"
$ws.Range("E46").Value = "-0.78%"
$ws.Range("E47").Value = "-0.77%"
$ws.Range("E48").Value = "20.89%"
$ws.Range("D49").Value = "0.003969"
$ws.Range("E49").Value = "-5.51%"
$ws.Range("D50").Value = "0.00002084"
$ws.Range("E50").Value = "-0.77%"
$ws.Range("D51").Value = "0.0001985"
$ws.Range("E51").Value = "-0.77%"

# Restore default (unstyled) cell formatting so output matches original styling.
$ws.Range("D2:E51").Style = "Normal"

